# "capitalization added to views" - tidy up the TODO list on Sheet1:
#  - clear the now-resolved "error checking for zero entries on web forms" item
#  - remove the three completed/obsolete TODO rows immediately below it
#    (the remaining rows shift up to close the gap)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the "error checking for zero entries on web forms" note, but keep
# the TODO header in D35 untouched.
$ws.Cells.Item(35, 5).ClearContents()

# Remove rows 36-38 entirely (need to add purchase date capture.../ need to
# sort capitalisation on words / fix order changes...) - the rows below
# shift up to fill the gap.
$ws.Rows("36:38").Delete()
